$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "test sheet 1" -> "test sheet 1 changed" (shared string used by A1)
$ws.Range("A1").Value = "test sheet 1 changed"

# Reset the saved selection/active cell back to A1 (was E6)
$ws.Range("A1").Select()
